$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: change a numeric/text cell's *value* while forcing its style (xf)
# to match a known-good template cell elsewhere on the sheet. This avoids the
# engine minting a brand-new style record (it reuses the template's existing
# style index) when a value's "natural" type doesn't match the cell's
# current number format (e.g. writing the literal text "0" into a cell that
# used to hold a real number, or vice versa).
# ---------------------------------------------------------------------------
function Set-CellWithStyle {
    param(
        [int]$Row,
        [int]$Col,
        $Value,
        [int]$StyleSrcRow,
        [int]$StyleSrcCol
    )
    $dst = $ws.Cells.Item($Row, $Col)
    if ($Value -is [string]) {
        # Force text interpretation so numeric-looking strings (e.g. "0")
        # are not silently coerced back into a number.
        $dst.NumberFormat = "@"
        $dst.Value = $Value
    } else {
        $dst.Value = $Value
    }
    $src = $ws.Cells.Item($StyleSrcRow, $StyleSrcCol)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# Template cells (row 15 is untouched by this edit) for each style used below
# s=14 (text, right aligned)      -> C15
# s=15 (number, #,##0)            -> I15
# s=16 (number, #,##0.0 w/ paren) -> K15

# ---------------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------------

# Police Commissioner name
$ws.Cells.Item(6, 13).Value = "Edward A. Caban"

# "Volume 30   Number  26" -> "...27"  (only the trailing issue number run)
$c8 = $ws.Cells.Item(8, 1)
$txt8 = $c8.Text
$start = $txt8.Length - 1
$c8.Characters($start, 2).Text = "27"

# "Report Covering the Week  6/26/2023  Through  7/2/2023"
#    -> "...7/3/2023  Through  7/9/2023"
$c9 = $ws.Cells.Item(9, 3)
$txt9 = $c9.Text
$idx1 = $txt9.IndexOf("6/26/2023") + 1
$c9.Characters($idx1, 9).Text = "7/3/2023"
$txt9b = $c9.Text
$idx2 = $txt9b.IndexOf("7/2/2023") + 1
$c9.Characters($idx2, 8).Text = "7/9/2023"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 6).Value = 4
$ws.Cells.Item(16, 8).Value = 33.333333333333
$ws.Cells.Item(16, 10).Value = 26
$ws.Cells.Item(16, 11).Value = 53.846153846153
$ws.Cells.Item(16, 12).Value = 53.846153846153
$ws.Cells.Item(16, 13).Value = -24.528301886792
$ws.Cells.Item(16, 14).Value = -86.885245901639

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 5).Value = -50
$ws.Cells.Item(17, 7).Value = 9
$ws.Cells.Item(17, 8).Value = 55.555555555555
$ws.Cells.Item(17, 9).Value = 80
$ws.Cells.Item(17, 10).Value = 51
$ws.Cells.Item(17, 11).Value = 56.862745098039
$ws.Cells.Item(17, 12).Value = 73.913043478260
$ws.Cells.Item(17, 13).Value = 42.857142857142
$ws.Cells.Item(17, 14).Value = -43.661971830985

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 100
$ws.Cells.Item(18, 7).Value = 8
$ws.Cells.Item(18, 8).Value = -37.5
$ws.Cells.Item(18, 9).Value = 58
$ws.Cells.Item(18, 10).Value = 41
$ws.Cells.Item(18, 11).Value = 41.463414634146
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -53.6
$ws.Cells.Item(18, 14).Value = -89.605734767025

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Cells.Item(19, 3).Value = 12
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 7).Value = 44
$ws.Cells.Item(19, 8).Value = -18.181818181818
$ws.Cells.Item(19, 9).Value = 230
$ws.Cells.Item(19, 10).Value = 290
$ws.Cells.Item(19, 11).Value = -20.689655172413
$ws.Cells.Item(19, 12).Value = 29.943502824858
$ws.Cells.Item(19, 13).Value = 45.569620253164
$ws.Cells.Item(19, 14).Value = -6.122448979591

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = -33.333333333333
$ws.Cells.Item(20, 6).Value = 10
$ws.Cells.Item(20, 8).Value = -16.666666666666
$ws.Cells.Item(20, 9).Value = 60
$ws.Cells.Item(20, 10).Value = 60
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 93.548387096774
$ws.Cells.Item(20, 13).Value = -28.571428571428
$ws.Cells.Item(20, 14).Value = -94.208494208494

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 18
$ws.Cells.Item(21, 4).Value = 21
$ws.Cells.Item(21, 5).Value = -14.285714285714
$ws.Cells.Item(21, 6).Value = 69
$ws.Cells.Item(21, 7).Value = 76
$ws.Cells.Item(21, 8).Value = -9.210526315789
$ws.Cells.Item(21, 9).Value = 476
$ws.Cells.Item(21, 10).Value = 476
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 37.971014492753
$ws.Cells.Item(21, 13).Value = -1.449275362318
$ws.Cells.Item(21, 14).Value = -79.331306990881

# ---------------------------------------------------------------------------
# Row 22 - Transit  (F22 flips from a number to the text "0")
# ---------------------------------------------------------------------------
Set-CellWithStyle -Row 22 -Col 6 -Value "0" -StyleSrcRow 15 -StyleSrcCol 3
$ws.Cells.Item(22, 13).Value = -20

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Cells.Item(24, 3).Value = 29
$ws.Cells.Item(24, 4).Value = 40
$ws.Cells.Item(24, 5).Value = -27.5
$ws.Cells.Item(24, 7).Value = 168
$ws.Cells.Item(24, 8).Value = -35.714285714285
$ws.Cells.Item(24, 9).Value = 792
$ws.Cells.Item(24, 10).Value = 951
$ws.Cells.Item(24, 11).Value = -16.719242902208
$ws.Cells.Item(24, 12).Value = 63.298969072165
$ws.Cells.Item(24, 13).Value = 31.780366056572

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 3
$ws.Cells.Item(25, 5).Value = -25
$ws.Cells.Item(25, 6).Value = 29
$ws.Cells.Item(25, 7).Value = 21
$ws.Cells.Item(25, 8).Value = 38.095238095238
$ws.Cells.Item(25, 9).Value = 199
$ws.Cells.Item(25, 10).Value = 161
$ws.Cells.Item(25, 11).Value = 23.602484472049
$ws.Cells.Item(25, 12).Value = 57.936507936507
$ws.Cells.Item(25, 13).Value = 5.851063829787

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*  (C,D flip number->text "0"; E flips number->text "***.*")
# ---------------------------------------------------------------------------
Set-CellWithStyle -Row 26 -Col 3 -Value "0" -StyleSrcRow 15 -StyleSrcCol 3
Set-CellWithStyle -Row 26 -Col 4 -Value "0" -StyleSrcRow 15 -StyleSrcCol 3
Set-CellWithStyle -Row 26 -Col 5 -Value "***.*" -StyleSrcRow 15 -StyleSrcCol 3

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes  (C flips text "0" -> number 1)
# ---------------------------------------------------------------------------
Set-CellWithStyle -Row 27 -Col 3 -Value 1 -StyleSrcRow 15 -StyleSrcCol 9
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 8).Value = 33.333333333333
$ws.Cells.Item(27, 9).Value = 14
$ws.Cells.Item(27, 10).Value = 24
$ws.Cells.Item(27, 11).Value = -41.666666666666
$ws.Cells.Item(27, 12).Value = 27.272727272727

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
#   C flips number 1 -> text "0"
#   D flips text "0" -> number 1
#   E flips text "***.*" -> number -100
#   G flips text "0" -> number 1
#   H flips text "***.*" -> number 200
# ---------------------------------------------------------------------------
Set-CellWithStyle -Row 28 -Col 3 -Value "0" -StyleSrcRow 15 -StyleSrcCol 3
Set-CellWithStyle -Row 28 -Col 4 -Value 1 -StyleSrcRow 15 -StyleSrcCol 9
Set-CellWithStyle -Row 28 -Col 5 -Value -100 -StyleSrcRow 15 -StyleSrcCol 11
Set-CellWithStyle -Row 28 -Col 7 -Value 1 -StyleSrcRow 15 -StyleSrcCol 9
Set-CellWithStyle -Row 28 -Col 8 -Value 200 -StyleSrcRow 15 -StyleSrcCol 11
$ws.Cells.Item(28, 10).Value = 3
$ws.Cells.Item(28, 11).Value = 0

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
#   C flips number 1 -> text "0"
#   D flips text "0" -> number 1
#   E flips text "***.*" -> number -100
#   G flips text "0" -> number 1
#   H flips text "***.*" -> number 100
# ---------------------------------------------------------------------------
Set-CellWithStyle -Row 29 -Col 3 -Value "0" -StyleSrcRow 15 -StyleSrcCol 3
Set-CellWithStyle -Row 29 -Col 4 -Value 1 -StyleSrcRow 15 -StyleSrcCol 9
Set-CellWithStyle -Row 29 -Col 5 -Value -100 -StyleSrcRow 15 -StyleSrcCol 11
Set-CellWithStyle -Row 29 -Col 7 -Value 1 -StyleSrcRow 15 -StyleSrcCol 9
Set-CellWithStyle -Row 29 -Col 8 -Value 100 -StyleSrcRow 15 -StyleSrcCol 11
$ws.Cells.Item(29, 10).Value = 3
$ws.Cells.Item(29, 11).Value = -33.333333333333
